$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 1483.4642
$ws.Range("I70").Value = 1067.4445
$ws.Range("J70").Value = 1680.5264
$ws.Range("K70").Value = 3202.3335
$ws.Range("L70").Value = 5041.5792
$ws.Range("M70").Value = -2932.3335
$ws.Range("N70").Value = -5581.5792
# Row 73
$ws.Range("H73").Value = 1483.4642
$ws.Range("I73").Value = 1067.4445
$ws.Range("J73").Value = 1680.5264
$ws.Range("K73").Value = 3202.3335
$ws.Range("L73").Value = 5041.5792
$ws.Range("M73").Value = -2266.3335
$ws.Range("N73").Value = -6913.5792
# Row 74
$ws.Range("H74").Value = 4208.727
$ws.Range("I74").Value = 3999.3333
$ws.Range("J74").Value = 4287.25
$ws.Range("K74").Value = 3999.3333
$ws.Range("L74").Value = 4287.25
$ws.Range("M74").Value = -3063.3333
$ws.Range("N74").Value = -6159.25
# Row 76
$ws.Range("H76").Value = 3054.5454
$ws.Range("J76").Value = 3114.2856
$ws.Range("L76").Value = 3114.2856
$ws.Range("N76").Value = -3744.2856
# Row 77
$ws.Range("H77").Value = 4208.727
$ws.Range("I77").Value = 3999.3333
$ws.Range("J77").Value = 4287.25
$ws.Range("K77").Value = 19996.6665
$ws.Range("L77").Value = 21436.25
$ws.Range("M77").Value = -15316.6665
$ws.Range("N77").Value = -30796.25
# Row 79
$ws.Range("H79").Value = 3054.5454
$ws.Range("J79").Value = 3114.2856
$ws.Range("L79").Value = 3114.2856
$ws.Range("N79").Value = -5298.2856
# Row 92
$ws.Range("H92").Value = 1043.0588
$ws.Range("I92").Value = 938
$ws.Range("J92").Value = 1533.3334
$ws.Range("K92").Value = 938
$ws.Range("L92").Value = 1533.3334
$ws.Range("M92").Value = 310
$ws.Range("N92").Value = -4029.3334

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 110
$ws.Range("H110").Value = 1626.4546
$ws.Range("I110").Value = 1599
$ws.Range("J110").Value = 1750
$ws.Range("K110").Value = 1599
$ws.Range("L110").Value = 1750
$ws.Range("M110").Value = 446
$ws.Range("N110").Value = -5840
# Row 122
$ws.Range("H122").Value = 2300
$ws.Range("I122").Value = 1885.7142
$ws.Range("J122").Value = 2714.2856
$ws.Range("K122").Value = 5657.142599999999
$ws.Range("L122").Value = 8142.8568
$ws.Range("M122").Value = -3207.142599999999
$ws.Range("N122").Value = -13042.8568
# Row 124
$ws.Range("H124").Value = 22145.7
$ws.Range("J124").Value = 22145.7
$ws.Range("L124").Value = 22145.7
$ws.Range("N124").Value = -31965.7
# Row 137
$ws.Range("H137").Value = 43571.285
$ws.Range("I137").Value = 39000
$ws.Range("J137").Value = 44333.168
$ws.Range("K137").Value = 39000
$ws.Range("L137").Value = 44333.168
$ws.Range("N137").Value = -54533.168
$ws.Range("M137").Value = -33900

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 57
$ws.Range("H57").Value = 47000
$ws.Range("J57").Value = 47000
$ws.Range("L57").Value = 47000
$ws.Range("N57").Value = -48440
# Row 81
$ws.Range("H81").Value = 25938.092
$ws.Range("J81").Value = 25938.092
$ws.Range("L81").Value = 25938.092
$ws.Range("N81").Value = -28060.092
# Row 84
$ws.Range("H84").Value = 25938.092
$ws.Range("J84").Value = 25938.092
$ws.Range("L84").Value = 77814.276
$ws.Range("N84").Value = -88422.276
# Row 135
$ws.Range("H135").Value = 30260
$ws.Range("J135").Value = 30260
$ws.Range("L135").Value = 30260
$ws.Range("N135").Value = -40400
# Row 136
$ws.Range("H136").Value = 47000
$ws.Range("J136").Value = 47000
$ws.Range("L136").Value = 47000
$ws.Range("N136").Value = -57200

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 62501212
$ws.Range("I31").Value = 71429290
$ws.Range("J31").Value = 50001910
$ws.Range("K31").Value = 71429290
$ws.Range("L31").Value = 50001910
$ws.Range("M31").Value = -71428995
$ws.Range("N31").Value = -50002500
# Row 34
$ws.Range("H34").Value = 62501212
$ws.Range("I34").Value = 71429290
$ws.Range("J34").Value = 50001910
$ws.Range("K34").Value = 71429290
$ws.Range("L34").Value = 50001910
$ws.Range("M34").Value = -71429088
$ws.Range("N34").Value = -50002314
# Row 124
$ws.Range("H124").Value = 25203.75
$ws.Range("J124").Value = 25203.75
$ws.Range("L124").Value = 25203.75
$ws.Range("N124").Value = -30113.75

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 92
$ws.Range("H92").Value = 1667424.5
$ws.Range("I92").Value = 498
$ws.Range("J92").Value = 2000809.8
$ws.Range("K92").Value = 1494
$ws.Range("L92").Value = 6002429.4
$ws.Range("M92").Value = -246
$ws.Range("N92").Value = -6004925.4
# Row 100
$ws.Range("H100").Value = 3775.3845
$ws.Range("J100").Value = 3775.3845
$ws.Range("L100").Value = 11326.1535
$ws.Range("N100").Value = -12948.1535
# Row 107
$ws.Range("H107").Value = 3333598.5
$ws.Range("I107").Value = 255
$ws.Range("J107").Value = 14444744
$ws.Range("K107").Value = 765
$ws.Range("L107").Value = 43334232
$ws.Range("M107").Value = 1155
$ws.Range("N107").Value = -43338072
# Row 108
$ws.Range("H108").Value = 1154.5
$ws.Range("I108").Value = 1154.5
$ws.Range("J108").Value = 0
$ws.Range("K108").Value = 3463.5
$ws.Range("L108").Value = 0
$ws.Range("M108").Value = -583.5
$ws.Range("N108").ClearContents()
# Row 109
$ws.Range("H109").Value = 2168.8333
$ws.Range("J109").Value = 2665.6
$ws.Range("L109").Value = 7996.799999999999
$ws.Range("N109").Value = -10076.8
# Row 110
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("M110").ClearContents()
# Row 111
$ws.Range("H111").Value = 17325.666
$ws.Range("I111").Value = 20430.8
$ws.Range("K111").Value = 61292.39999999999
$ws.Range("M111").Value = -58225.39999999999
# Row 112
$ws.Range("H112").Value = 66670720
$ws.Range("I112").Value = 2000
$ws.Range("J112").Value = 71432770
$ws.Range("K112").Value = 6000
$ws.Range("L112").Value = 214298310
$ws.Range("M112").Value = -4892
$ws.Range("N112").Value = -214300526

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 51
$ws.Range("H51").Value = 25160
$ws.Range("J51").Value = 25160
$ws.Range("L51").Value = 25160
$ws.Range("N51").Value = -26178
# Row 80
$ws.Range("H80").Value = 3113.1333
$ws.Range("I80").Value = 2600
$ws.Range("J80").Value = 3882.8333
$ws.Range("K80").Value = 2600
$ws.Range("L80").Value = 3882.8333
$ws.Range("M80").Value = -1602
$ws.Range("N80").Value = -5878.8333
# Row 83
$ws.Range("H83").Value = 3113.1333
$ws.Range("I83").Value = 2600
$ws.Range("J83").Value = 3882.8333
$ws.Range("K83").Value = 13000
$ws.Range("L83").Value = 19414.1665
$ws.Range("M83").Value = -8008
$ws.Range("N83").Value = -29398.1665
# Row 132
$ws.Range("H132").Value = 2461.5386
$ws.Range("I132").Value = 2182.8857
$ws.Range("K132").Value = 6548.657099999999
$ws.Range("M132").Value = -4018.657099999999

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 68
$ws.Range("H68").Value = 12154.546
$ws.Range("I68").Value = 13744.444
$ws.Range("J68").Value = 5000
$ws.Range("K68").Value = 13744.444
$ws.Range("L68").Value = 5000
$ws.Range("M68").Value = -12995.444
$ws.Range("N68").Value = -6498
# Row 71
$ws.Range("H71").Value = 12154.546
$ws.Range("I71").Value = 13744.444
$ws.Range("J71").Value = 5000
$ws.Range("K71").Value = 68722.22
$ws.Range("L71").Value = 25000
$ws.Range("M71").Value = -64978.22
$ws.Range("N71").Value = -32488
# Row 122
$ws.Range("H122").Value = 6457.2856
$ws.Range("I122").Value = 6866.3335
$ws.Range("J122").Value = 5721
$ws.Range("K122").Value = 20599.0005
$ws.Range("L122").Value = 17163
$ws.Range("M122").Value = -18149.0005
$ws.Range("N122").Value = -22063
# Row 127
$ws.Range("H127").Value = 33269.855
$ws.Range("J127").Value = 33269.855
$ws.Range("L127").Value = 33269.855
$ws.Range("N127").Value = -43189.855

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 22838546
$ws.Range("I122").Value = 83335336
$ws.Range("J122").Value = 152251.62
$ws.Range("K122").Value = 250006008
$ws.Range("L122").Value = 456754.86
$ws.Range("M122").Value = -250003558
$ws.Range("N122").Value = -461654.86

Write-Host "Applied scheduled market data update."